$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.8243956666666667
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.833274
$ws.Range("N2").Value = 2.499822
$ws.Range("O2").Value = 0.8541534975292705
$ws.Range("P2").Value = 0.8541534975292705
$ws.Range("Q2").Value = 0.686947474746
$ws.Range("R2").Value = 6.182527272714
$ws.Range("S2").Value = 0.8541534975292705
$ws.Range("T2").Value = 0.8541534975292705

# Row 3 updates
$ws.Range("G3").Value = 0.8243956666666667
$ws.Range("M3").Value = 0.1422813333333333
$ws.Range("N3").Value = 0.426844
$ws.Range("O3").Value = 0.1458465024707295
$ws.Range("P3").Value = 0.1458465024707295
$ws.Range("Q3").Value = 0.1172961146475556
$ws.Range("R3").Value = 1.055665031828
$ws.Range("S3").Value = 0.1458465024707295
$ws.Range("T3").Value = 0.1458465024707295
